$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: extend with two more sequential values (I1=9, J1=10)
$ws.Range("I1").Value = 9
$ws.Range("J1").Value = 10

# Row 2: update existing values to new "random limit" values and add two more columns
$ws.Range("A2").Value = 577
$ws.Range("B2").Value = 445
$ws.Range("C2").Value = 460
$ws.Range("D2").Value = 605
$ws.Range("E2").Value = 458
$ws.Range("F2").Value = 500
$ws.Range("G2").Value = 492
$ws.Range("H2").Value = 485
$ws.Range("I2").Value = 589
$ws.Range("J2").Value = 501
